# BCA_marks.xlsx update
# - Column G (SGPA?) values on Sheet1, rows 2-76, are each increased by 0.154
# - The active selection on the sheet is changed to the full column H (H1:H1048576),
#   with H1 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 76; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = $cell.Value() + 0.154
}

$ws.Range("H1:H1048576").Select()
